$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '90.436.46'
$ws.Range("E2").Value = '  +1.03%  '

# Row 3
$ws.Range("D3").Value = '3.044.68'
$ws.Range("E3").Value = '  -3.20%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.85%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.72%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '611.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.369'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.59%  '

# Row 8
$ws.Range("E8").Value = '  +11.97%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.21%  '

# Row 10
$ws.Range("D10").Value = '3.043.87'
$ws.Range("E10").Value = '  -3.10%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.676'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +19.67%  '

# Row 12
$ws.Range("E12").Value = '  +5.01%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000242'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.92%  '

# Row 14
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '90.280.00'
$ws.Range("E14").Value = '  +1.18%  '

# Row 15
$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.33'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.33%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.35%  '

# Row 17
$ws.Range("D17").Value = '3.614.32'
$ws.Range("E17").Value = '  -2.54%  '

# Row 18
$ws.Range("D18").Value = '2.997.58'
$ws.Range("E18").Value = '  -4.49%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.94%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000221'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.16%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.42%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '426.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.40%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.60%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.52%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.63%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '83.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.54%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.67%  '

# Row 28
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '3.206.11'
$ws.Range("E28").Value = '  -2.78%  '

# Row 29
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.04%  '

# Row 30
$ws.Range("E30").Value = '  +3.38%  '

# Row 31
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.47%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.86%  '

# Row 33
$ws.Range("B33").Value = 'dogwifhat'
$ws.Range("C33").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.01%  '

# Row 34
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '503.23'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.52%  '

# Row 35
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.68'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.22%  '

# Row 36
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.81'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.61%  '

# Row 37
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.66%  '

# Row 38
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.25'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.79%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.133'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -10.43%  '

# Row 40
$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.10%  '

# Row 41
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("E41").Value = '  +0.34%  '

# Row 42
$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.01%  '

# Row 43
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.137'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.37%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.360'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.11%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.55%  '

# Row 46
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '143.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.84%  '

# Row 47
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0697'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.10%  '

# Row 48
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '43.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.26%  '

# Row 49
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.38%  '

# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '160.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.06%  '

# Row 51
$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.32%  '
